# SYM_MODEL.xlsx update script
# Applies Q3'2025 actuals, new "Beat Guidance" notes, a new "EV/Backlog" metric,
# and refreshes the active sheet / selection state to match the authored edit.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1. "Main" sheet - refreshed price / share count / debt figures
# ----------------------------------------------------------------------
$main = $wb.Worksheets.Item("Main")

$main.Range("B4").Value = 52          # Price
$main.Range("G4").Value = 109.2       # Class A Common market cap helper
$main.Range("C5").Value = "Q32025"    # Shares outstanding "as of" quarter
$main.Range("B7").Value = 778         # Cash + PPE
# B6 (Market Cap = B4*B5) and B9 (Enterprise Value) are formulas and
# recompute automatically.

# ----------------------------------------------------------------------
# 2. "Model" sheet - add the Q3 2025 (column R) actuals
# ----------------------------------------------------------------------
$model = $wb.Worksheets.Item("Model")

# Revenue build
$model.Range("R3").Value = 559        # Systems
$model.Range("R4").Value = 8          # Software Maintenance and support
$model.Range("R5").Value = 25         # Operation Services
$model.Range("R7").Formula = "=SUM(R3:R5)"   # Revenue

# Cost build
$model.Range("R8").Value = 458        # Systems Cost
$model.Range("R9").Value = 2          # Software Cost
$model.Range("R10").Value = 24        # Operations Cost
$model.Range("R11").Formula = "=R7 - SUM(R8:R10)"   # Gross Profit

# Opex
$model.Range("R12").Value = 52        # R&D
$model.Range("R13").Value = 92        # SGA
$model.Range("R14").Formula = "=R11-SUM(R12:R13)"   # Operating Income

# Below the line
$model.Range("R15").Value = -28       # Income Before taxes
$model.Range("R16").Value = 0.04      # Tax
$model.Range("R17").Value = 4         # Loss From Equity Method Investment
$model.Range("R18").Formula = "=R15-R16-R17"         # Net Income

$model.Range("R19").Value = -0.05     # Basic EPS
$model.Range("R20").Value = -0.05     # Diluted EPS

# Growth metrics
$model.Range("R22").Formula = "=(R7/N7) - 1"   # Revenue Growth YOY Q
$model.Range("S22").Formula = "=(S7/O7) - 1"
$model.Range("R23").Formula = "= (R7/Q7) - 1"  # Revenue Growth last Q
$model.Range("S23").Formula = "= (S7/P7) - 1"

# Margins
$model.Range("R25").Value = 0.182     # Gross Margin
$model.Range("R26").Value = 0.215     # Adjusted Gross Margin
$model.Range("R27").Value = 45        # Adjusted EBITDA

# Cash flow
$model.Range("R29").Value = -138      # Cash Flow from Operations
$model.Range("R30").Value = 15        # Capital Expendature
$model.Range("R31").Formula = "=R29-R30"             # FCF

# Highlight the newly reported revenue figure (bold + italic, like the
# author's "beat" callout)
$model.Range("R7").Font.Bold = $true
$model.Range("R7").Font.Italic = $true

# Guidance notes
$model.Range("Q36").Value = "Beat Guidance"
$model.Range("R36").Value = "590 to 610 mill"
$model.Range("R37").Value = "adjusted ebitda of 45-49"

# ----------------------------------------------------------------------
# 3. "Valuation Metrics" sheet - add EV/Backlog metric
# ----------------------------------------------------------------------
$valMetrics = $wb.Worksheets.Item("Valuation Metrics")
$valMetrics.Range("A8").Value = "EV/Backlog"
$valMetrics.Range("B8").Value = 1.4

# ----------------------------------------------------------------------
# 4. View / selection state
# ----------------------------------------------------------------------
# Model sheet: scroll over and select R31
$model.Activate()
$excel.ActiveWindow.ScrollColumn = 10
$model.Range("R31").Select()

# Valuation Metrics: leave selection on the new last row, but it is no
# longer the active tab
$valMetrics.Range("A9").Select()

# Main becomes the active tab, with C6 selected
$main.Activate()
$main.Range("C6").Select()

Write-Output "Edit applied"
